$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style re-assignments (table "theme" GUIDs swapped to new values)
# ---------------------------------------------------------------------------
# Slides whose table currently uses style {6ED22EAA-...} -> {52057CED-...}
foreach ($idx in 4, 9, 10, 12, 19, 20) {
    $s = $p.Slides.Item($idx)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle("{52057CED-BF47-4F4C-8308-762A83911CC7}")
        }
    }
}

# Slides whose table currently uses style {81A5FA6C-...} -> {3F4F164F-...}
foreach ($idx in 6, 7, 13, 16) {
    $s = $p.Slides.Item($idx)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle("{3F4F164F-118B-4B08-8392-745AF1C78D78}")
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 12: fix the quoted error percentages
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$s12.Shapes.Item(1).TextFrame.TextRange.Text = "The percentages of 'cant_solve' and 'corrupt_data' to happen are 0.0187% and 0.0044% respectively."

# ---------------------------------------------------------------------------
# 3) Slide 15: re-word the title ("Task 4 - Annotators' results assessment"
#    -> "Task 4 - Annotator average result assessment") while keeping the
#    original two-run split point shifted to right after "Annotator".
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$titleShape = $s15.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

# First run currently spans characters 1-9 ("Task 4 - ")
$run1 = $titleRange.Characters(1, 9)
$run1.Text = "Task 4 - Annotator"

# After growing run1 to 19 characters, the remainder of the original second
# run ("Annotators' results assessment", 31 chars) now starts at character 20
$titleRange2 = $titleShape.TextFrame.TextRange
$run2 = $titleRange2.Characters(20, 31)
$run2.Text = " average result assessment"

# ---------------------------------------------------------------------------
# 4) Slide 15: adjust the cropped screenshot picture (re-crop + reposition)
# ---------------------------------------------------------------------------
$pic = $s15.Shapes.Item(4)
$pf = $pic.PictureFormat
$pf.CropTop = 58.8348
$pf.CropBottom = 57.79935

$pic.Left = (738525.0 / 12700.0) + 0.00002
$pic.Top = (1302825.0 / 12700.0) + 0.00002
$pic.Width = (3152824.0 / 12700.0) + 0.00002
$pic.Height = (3266050.0 / 12700.0) + 0.00002
